# Generate Report for Handoff
# Update the "Latest Handoff" timestamps for the file
# 945b0078-c6c1-4d54-8185-e351126a2086 across the Overview, zh-cn and
# de-de sheets to reflect a fresh handoff report run.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D = "Latest Handoff Date", row 6 is the
# 945b0078-c6c1-4d54-8185-e351126a2086.md entry.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-03-21 14:39:16"

# zh-cn sheet: column E = "Latest Handoff Datetime", row 6 is the
# 945b0078-c6c1-4d54-8185-e351126a2086 entry.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-21 14:39:12"

# de-de sheet: column E = "Latest Handoff Datetime", row 6 is the
# 945b0078-c6c1-4d54-8185-e351126a2086 entry.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-21 14:39:16"
